# Update AdultChildFlag / PSA (country & gender dependent) and RootSearch2
# employment-alignment rates across the employment_* worksheets.
#
# - employment_smales / employment_sfemales / employment_couples:
#     trim series down to 2011-2023 (was 2010-2035), recompute values,
#     relabel header to "empl_share", and employment_smales becomes the
#     active tab/selection.
# - employment_acfemales / employment_acmales / employment_femalewdep /
#     employment_malewdep: recompute 2011-2023 values, drop the explicit
#     "Normal 2" cell style, and normalise the sheet selection to A1:B14.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# employment_smales (sheet7): 2010-2035 -> 2011-2023, new empl_share data
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("employment_smales")
$ws.Range("A15:B27").ClearContents()
$ws.Range("B1").Value = "empl_share"
$ws.Range("A2").Value = 2011
$ws.Range("B2").Value = 0.44670278305339939
$ws.Range("A3").Value = 2012
$ws.Range("B3").Value = 0.44336450393902516
$ws.Range("A4").Value = 2013
$ws.Range("B4").Value = 0.44440277662332139
$ws.Range("A5").Value = 2014
$ws.Range("B5").Value = 0.45769354631775144
$ws.Range("A6").Value = 2015
$ws.Range("B6").Value = 0.46857670454141603
$ws.Range("A7").Value = 2016
$ws.Range("B7").Value = 0.46135340171417216
$ws.Range("A8").Value = 2017
$ws.Range("B8").Value = 0.48381900646628939
$ws.Range("A9").Value = 2018
$ws.Range("B9").Value = 0.4911795212142166
$ws.Range("A10").Value = 2019
$ws.Range("B10").Value = 0.48726882873224336
$ws.Range("A11").Value = 2020
$ws.Range("B11").Value = 0.49614597341437744
$ws.Range("A12").Value = 2021
$ws.Range("B12").Value = 0.52459623426267399
$ws.Range("A13").Value = 2022
$ws.Range("B13").Value = 0.53915577006649817
$ws.Range("A14").Value = 2023
$ws.Range("B14").Value = 0.52156811285163618

# ---------------------------------------------------------------------
# employment_sfemales (sheet8): 2010-2035 -> 2011-2023, new empl_share data
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("employment_sfemales")
$ws.Range("A15:B27").ClearContents()
$ws.Range("B1").Value = "empl_share"
$ws.Range("A2").Value = 2011
$ws.Range("B2").Value = 0.31394928614821893
$ws.Range("A3").Value = 2012
$ws.Range("B3").Value = 0.30848042060514502
$ws.Range("A4").Value = 2013
$ws.Range("B4").Value = 0.30956841493744147
$ws.Range("A5").Value = 2014
$ws.Range("B5").Value = 0.31944256894108353
$ws.Range("A6").Value = 2015
$ws.Range("B6").Value = 0.31710955870099683
$ws.Range("A7").Value = 2016
$ws.Range("B7").Value = 0.31537141978390892
$ws.Range("A8").Value = 2017
$ws.Range("B8").Value = 0.31503917525773195
$ws.Range("A9").Value = 2018
$ws.Range("B9").Value = 0.31132235685019893
$ws.Range("A10").Value = 2019
$ws.Range("B10").Value = 0.3186277678079098
$ws.Range("A11").Value = 2020
$ws.Range("B11").Value = 0.31556614940666167
$ws.Range("A12").Value = 2021
$ws.Range("B12").Value = 0.32083228970536581
$ws.Range("A13").Value = 2022
$ws.Range("B13").Value = 0.34580404255112251
$ws.Range("A14").Value = 2023
$ws.Range("B14").Value = 0.34735164074624669

# ---------------------------------------------------------------------
# employment_couples (sheet9): 2010-2035 -> 2011-2023, new empl_share data
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("employment_couples")
$ws.Range("A15:B27").ClearContents()
$ws.Range("B1").Value = "empl_share"
$ws.Range("A2").Value = 2011
$ws.Range("B2").Value = 0.97541612386703491
$ws.Range("A3").Value = 2012
$ws.Range("B3").Value = 0.97455012798309326
$ws.Range("A4").Value = 2013
$ws.Range("B4").Value = 0.96620428562164307
$ws.Range("A5").Value = 2014
$ws.Range("B5").Value = 0.97489017248153687
$ws.Range("A6").Value = 2015
$ws.Range("B6").Value = 0.9752197265625
$ws.Range("A7").Value = 2016
$ws.Range("B7").Value = 0.98048520088195801
$ws.Range("A8").Value = 2017
$ws.Range("B8").Value = 0.97633326053619385
$ws.Range("A9").Value = 2018
$ws.Range("B9").Value = 0.98472893238067627
$ws.Range("A10").Value = 2019
$ws.Range("B10").Value = 0.98302763700485229
$ws.Range("A11").Value = 2020
$ws.Range("B11").Value = 0.98899710178375244
$ws.Range("A12").Value = 2021
$ws.Range("B12").Value = 0.99124062061309814
$ws.Range("A13").Value = 2022
$ws.Range("B13").Value = 0.9915042519569397
$ws.Range("A14").Value = 2023
$ws.Range("B14").Value = 0.99202412366867065

# ---------------------------------------------------------------------
# employment_acfemales (sheet10): recomputed values, drop explicit style
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("employment_acfemales")
$ws.Range("A2").Value = 2011
$ws.Range("B2").Value = 0.36003186935365528
$ws.Range("A3").Value = 2012
$ws.Range("B3").Value = 0.38692491437698423
$ws.Range("A4").Value = 2013
$ws.Range("B4").Value = 0.36132873647659114
$ws.Range("A5").Value = 2014
$ws.Range("B5").Value = 0.38088347179059528
$ws.Range("A6").Value = 2015
$ws.Range("B6").Value = 0.38498832940988709
$ws.Range("A7").Value = 2016
$ws.Range("B7").Value = 0.4058379925777626
$ws.Range("A8").Value = 2017
$ws.Range("B8").Value = 0.43108781433921883
$ws.Range("A9").Value = 2018
$ws.Range("B9").Value = 0.40755919614401981
$ws.Range("A10").Value = 2019
$ws.Range("B10").Value = 0.42056689824920096
$ws.Range("A11").Value = 2020
$ws.Range("B11").Value = 0.45737911395870984
$ws.Range("A12").Value = 2021
$ws.Range("B12").Value = 0.45357219704792257
$ws.Range("A13").Value = 2022
$ws.Range("B13").Value = 0.4353627057118411
$ws.Range("A14").Value = 2023
$ws.Range("B14").Value = 0.41823639194988294
$ws.Range("A1:B14").Style = "Normal"

# ---------------------------------------------------------------------
# employment_acmales (sheet11): recomputed values, drop explicit style
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("employment_acmales")
$ws.Range("A2").Value = 2011
$ws.Range("B2").Value = 0.51128277630011854
$ws.Range("A3").Value = 2012
$ws.Range("B3").Value = 0.51393263199108286
$ws.Range("A4").Value = 2013
$ws.Range("B4").Value = 0.5076697666645027
$ws.Range("A5").Value = 2014
$ws.Range("B5").Value = 0.51067476163189951
$ws.Range("A6").Value = 2015
$ws.Range("B6").Value = 0.53577002132648832
$ws.Range("A7").Value = 2016
$ws.Range("B7").Value = 0.58570650522171619
$ws.Range("A8").Value = 2017
$ws.Range("B8").Value = 0.59110843342391406
$ws.Range("A9").Value = 2018
$ws.Range("B9").Value = 0.56825257026071563
$ws.Range("A10").Value = 2019
$ws.Range("B10").Value = 0.57190222490428488
$ws.Range("A11").Value = 2020
$ws.Range("B11").Value = 0.56928754219274136
$ws.Range("A12").Value = 2021
$ws.Range("B12").Value = 0.61268070477504499
$ws.Range("A13").Value = 2022
$ws.Range("B13").Value = 0.59899189733983293
$ws.Range("A14").Value = 2023
$ws.Range("B14").Value = 0.59106341164084708
$ws.Range("A1:B14").Style = "Normal"

# ---------------------------------------------------------------------
# employment_femalewdep (sheet12): recomputed values, drop explicit style
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("employment_femalewdep")
$ws.Range("A2").Value = 2011
$ws.Range("B2").Value = 0.61680769920349121
$ws.Range("A3").Value = 2012
$ws.Range("B3").Value = 0.61431878805160522
$ws.Range("A4").Value = 2013
$ws.Range("B4").Value = 0.62226575613021851
$ws.Range("A5").Value = 2014
$ws.Range("B5").Value = 0.56901419162750244
$ws.Range("A6").Value = 2015
$ws.Range("B6").Value = 0.55749976634979248
$ws.Range("A7").Value = 2016
$ws.Range("B7").Value = 0.61379510164260864
$ws.Range("A8").Value = 2017
$ws.Range("B8").Value = 0.60692942142486572
$ws.Range("A9").Value = 2018
$ws.Range("B9").Value = 0.62670052051544189
$ws.Range("A10").Value = 2019
$ws.Range("B10").Value = 0.70763528347015381
$ws.Range("A11").Value = 2020
$ws.Range("B11").Value = 0.69641482830047607
$ws.Range("A12").Value = 2021
$ws.Range("B12").Value = 0.72090119123458862
$ws.Range("A13").Value = 2022
$ws.Range("B13").Value = 0.72393757104873657
$ws.Range("A14").Value = 2023
$ws.Range("B14").Value = 0.73154628276824951
$ws.Range("A1:B14").Style = "Normal"

# ---------------------------------------------------------------------
# employment_malewdep (sheet13): recomputed values, drop explicit style
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("employment_malewdep")
$ws.Range("A2").Value = 2011
$ws.Range("B2").Value = 0.83789324760437012
$ws.Range("A3").Value = 2012
$ws.Range("B3").Value = 0.80550128221511841
$ws.Range("A4").Value = 2013
$ws.Range("B4").Value = 0.7992321252822876
$ws.Range("A5").Value = 2014
$ws.Range("B5").Value = 0.79625046253204346
$ws.Range("A6").Value = 2015
$ws.Range("B6").Value = 0.81983458995819092
$ws.Range("A7").Value = 2016
$ws.Range("B7").Value = 0.81091392040252686
$ws.Range("A8").Value = 2017
$ws.Range("B8").Value = 0.82456725835800171
$ws.Range("A9").Value = 2018
$ws.Range("B9").Value = 0.85373425483703613
$ws.Range("A10").Value = 2019
$ws.Range("B10").Value = 0.84969794750213623
$ws.Range("A11").Value = 2020
$ws.Range("B11").Value = 0.82026958465576172
$ws.Range("A12").Value = 2021
$ws.Range("B12").Value = 0.86517715454101562
$ws.Range("A13").Value = 2022
$ws.Range("B13").Value = 0.88653171062469482
$ws.Range("A14").Value = 2023
$ws.Range("B14").Value = 0.92414480447769165
$ws.Range("A1:B14").Style = "Normal"

# ---------------------------------------------------------------------
# Sheet selections: every employment_* sheet now shows A1:B14 selected;
# employment_smales becomes the active tab (was employment_malewdep).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("employment_acfemales").Range("A1:B14").Select()
$wb.Worksheets.Item("employment_acmales").Range("A1:B14").Select()
$wb.Worksheets.Item("employment_femalewdep").Range("A1:B14").Select()
$wb.Worksheets.Item("employment_malewdep").Range("A1:B14").Select()
$wb.Worksheets.Item("employment_sfemales").Range("A1:B14").Select()
$wb.Worksheets.Item("employment_couples").Range("A1:B14").Select()

$wsActive = $wb.Worksheets.Item("employment_smales")
$wsActive.Range("A1:B14").Select()
$wsActive.Activate()
